$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated odds values per row, matching the target OOXML diff.
# Columns are referenced by their 1-based index (F=6, G=7, ... AO=41).

# Row 2
$ws.Cells.Item(2, 6).Value = 4.3  # F2
$ws.Cells.Item(2, 7).Value = 6.2  # G2
$ws.Cells.Item(2, 8).Value = 1.71  # H2
$ws.Cells.Item(2, 9).Value = 1.91  # I2
$ws.Cells.Item(2, 12).Value = 1.31  # L2
$ws.Cells.Item(2, 13).Value = 1.04  # M2
$ws.Cells.Item(2, 14).Value = 3.6  # N2
$ws.Cells.Item(2, 17).Value = 1.69  # Q2
$ws.Cells.Item(2, 18).Value = 1.36  # R2
$ws.Cells.Item(2, 20).Value = 1.76  # T2
$ws.Cells.Item(2, 21).Value = 2.12  # U2
$ws.Cells.Item(2, 22).Value = 2.1  # V2
$ws.Cells.Item(2, 23).Value = 1.19  # W2
$ws.Cells.Item(2, 24).Value = 21  # X2
$ws.Cells.Item(2, 25).Value = 11.5  # Y2
$ws.Cells.Item(2, 26).Value = 13.5  # Z2
$ws.Cells.Item(2, 27).Value = 22  # AA2
$ws.Cells.Item(2, 28).Value = 23  # AB2
$ws.Cells.Item(2, 29).Value = 11  # AC2
$ws.Cells.Item(2, 30).Value = 12  # AD2
$ws.Cells.Item(2, 31).Value = 21  # AE2
$ws.Cells.Item(2, 32).Value = 48  # AF2
$ws.Cells.Item(2, 33).Value = 24  # AG2
$ws.Cells.Item(2, 34).Value = 22  # AH2
$ws.Cells.Item(2, 35).Value = 38  # AI2

# Row 3
$ws.Cells.Item(3, 7).Value = 2.92  # G3
$ws.Cells.Item(3, 12).Value = 1.01  # L3
$ws.Cells.Item(3, 13).Value = 1.14  # M3
$ws.Cells.Item(3, 14).Value = 2.24  # N3
$ws.Cells.Item(3, 15).Value = 1.69  # O3
$ws.Cells.Item(3, 18).Value = 1.14  # R3
$ws.Cells.Item(3, 19).Value = 6.8  # S3
$ws.Cells.Item(3, 20).Value = 2.3  # T3
$ws.Cells.Item(3, 21).Value = 1.64  # U3
$ws.Cells.Item(3, 22).Value = 1.39  # V3
$ws.Cells.Item(3, 23).Value = 1.52  # W3
$ws.Cells.Item(3, 24).Value = 7  # X3
$ws.Cells.Item(3, 25).Value = 9.6  # Y3
$ws.Cells.Item(3, 26).Value = 25  # Z3
$ws.Cells.Item(3, 27).Value = 1000  # AA3
$ws.Cells.Item(3, 28).Value = 8.199999999999999  # AB3
$ws.Cells.Item(3, 29).Value = 8  # AC3
$ws.Cells.Item(3, 30).Value = 19.5  # AD3
$ws.Cells.Item(3, 31).Value = 1000  # AE3
$ws.Cells.Item(3, 32).Value = 18.5  # AF3
$ws.Cells.Item(3, 33).Value = 16.5  # AG3
$ws.Cells.Item(3, 34).Value = 32  # AH3
$ws.Cells.Item(3, 35).Value = 120  # AI3
$ws.Cells.Item(3, 36).Value = 60  # AJ3
$ws.Cells.Item(3, 37).Value = 60  # AK3
$ws.Cells.Item(3, 38).Value = 110  # AL3
$ws.Cells.Item(3, 39).Value = 310  # AM3
$ws.Cells.Item(3, 40).Value = 1000  # AN3
$ws.Cells.Item(3, 41).Value = 1000  # AO3

# Row 4
$ws.Cells.Item(4, 6).Value = 1.66  # F4
$ws.Cells.Item(4, 7).Value = 1.82  # G4
$ws.Cells.Item(4, 9).Value = 7  # I4
$ws.Cells.Item(4, 10).Value = 3.8  # J4
$ws.Cells.Item(4, 12).Value = 1.28  # L4
$ws.Cells.Item(4, 13).Value = 1.06  # M4
$ws.Cells.Item(4, 14).Value = 3.35  # N4
$ws.Cells.Item(4, 15).Value = 1.23  # O4
$ws.Cells.Item(4, 18).Value = 1.32  # R4
$ws.Cells.Item(4, 19).Value = 2.78  # S4
$ws.Cells.Item(4, 20).Value = 1.61  # T4
$ws.Cells.Item(4, 21).Value = 1.72  # U4
$ws.Cells.Item(4, 22).Value = 1.18  # V4
$ws.Cells.Item(4, 23).Value = 2.2  # W4
$ws.Cells.Item(4, 24).Value = 1000  # X4
$ws.Cells.Item(4, 25).Value = 28  # Y4
$ws.Cells.Item(4, 26).Value = 1000  # Z4
$ws.Cells.Item(4, 27).Value = 1000  # AA4
$ws.Cells.Item(4, 28).Value = 12  # AB4
$ws.Cells.Item(4, 29).Value = 13  # AC4
$ws.Cells.Item(4, 30).Value = 30  # AD4
$ws.Cells.Item(4, 31).Value = 1000  # AE4
$ws.Cells.Item(4, 32).Value = 14.5  # AF4
$ws.Cells.Item(4, 33).Value = 13.5  # AG4
$ws.Cells.Item(4, 34).Value = 29  # AH4
$ws.Cells.Item(4, 35).Value = 1000  # AI4
$ws.Cells.Item(4, 36).Value = 24  # AJ4
$ws.Cells.Item(4, 37).Value = 25  # AK4
$ws.Cells.Item(4, 38).Value = 1000  # AL4
$ws.Cells.Item(4, 39).Value = 1000  # AM4
$ws.Cells.Item(4, 40).Value = 1000  # AN4
$ws.Cells.Item(4, 41).Value = 1000  # AO4

# Row 5
$ws.Cells.Item(5, 7).Value = 1.85  # G5
$ws.Cells.Item(5, 10).Value = 3.6  # J5
$ws.Cells.Item(5, 12).Value = 1.01  # L5
$ws.Cells.Item(5, 14).Value = 3.05  # N5
$ws.Cells.Item(5, 15).Value = 1.35  # O5
$ws.Cells.Item(5, 16).Value = 1.79  # P5
$ws.Cells.Item(5, 18).Value = 1.27  # R5
$ws.Cells.Item(5, 19).Value = 3.7  # S5
$ws.Cells.Item(5, 20).Value = 1.94  # T5
$ws.Cells.Item(5, 22).Value = 1.2  # V5
$ws.Cells.Item(5, 23).Value = 2.18  # W5
$ws.Cells.Item(5, 25).Value = 1000  # Y5
$ws.Cells.Item(5, 28).Value = 8.6  # AB5

# Row 8
$ws.Cells.Item(8, 17).Value = 2.92  # Q8

# Row 9
$ws.Cells.Item(9, 9).Value = 870  # I9
$ws.Cells.Item(9, 14).Value = 2.54  # N9
$ws.Cells.Item(9, 15).Value = 1.48  # O9
$ws.Cells.Item(9, 17).Value = 2.18  # Q9

# Row 11
$ws.Cells.Item(11, 7).Value = 4.5  # G11
$ws.Cells.Item(11, 9).Value = 2.66  # I11
$ws.Cells.Item(11, 10).Value = 2.36  # J11

# Row 12
$ws.Cells.Item(12, 6).Value = 2.22  # F12
$ws.Cells.Item(12, 7).Value = 2.32  # G12
$ws.Cells.Item(12, 8).Value = 3.65  # H12
$ws.Cells.Item(12, 11).Value = 3.55  # K12
$ws.Cells.Item(12, 16).Value = 1.81  # P12

# Row 14
$ws.Cells.Item(14, 6).Value = 1.64  # F14
$ws.Cells.Item(14, 7).Value = 1.67  # G14
$ws.Cells.Item(14, 9).Value = 6.8  # I14
$ws.Cells.Item(14, 10).Value = 3.9  # J14
$ws.Cells.Item(14, 11).Value = 4.3  # K14
